# Auto-generated Excel COM-interop script to apply scheduled-runner market data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1700
$ws.Range("J43").Value = 1380
$ws.Range("L43").Value = 1380
$ws.Range("N43").Value = -1518
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("H132").Value = 13353.314
$ws.Range("I132").Value = 2105.1875
$ws.Range("K132").Value = 6315.5625
$ws.Range("M132").Value = -3785.5625
$ws.Range("H135").Value = 716751.4
$ws.Range("I135").Value = 835209.9399999999
$ws.Range("J135").Value = 6000
$ws.Range("K135").Value = 7516889.459999999
$ws.Range("L135").Value = 54000
$ws.Range("M135").Value = -7514354.459999999
$ws.Range("N135").Value = -59070
$ws.Range("H137").Value = 3886.853
$ws.Range("I137").Value = 4055.6072
$ws.Range("K137").Value = 12166.8216
$ws.Range("M137").Value = -9616.821599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2797.2222
$ws.Range("I61").Value = 2797.2222
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2797.2222
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2585.2222
$ws.Range("N61").Value = $null
$ws.Range("H74").Value = 1662.95
$ws.Range("I74").Value = 1554.2142
$ws.Range("K74").Value = 1554.2142
$ws.Range("M74").Value = -680.2141999999999
$ws.Range("H77").Value = 1662.95
$ws.Range("I77").Value = 1554.2142
$ws.Range("K77").Value = 7771.071
$ws.Range("M77").Value = -3403.071
$ws.Range("H122").Value = 6863
$ws.Range("I122").Value = 8799.875
$ws.Range("J122").Value = 5830
$ws.Range("K122").Value = 26399.625
$ws.Range("L122").Value = 17490
$ws.Range("M122").Value = -23949.625
$ws.Range("N122").Value = -22390
$ws.Range("H132").Value = 3820.5417
$ws.Range("I132").Value = 3531.3901
$ws.Range("K132").Value = 10594.1703
$ws.Range("M132").Value = -8064.1703
$ws.Range("H136").Value = 2797.2222
$ws.Range("I136").Value = 2797.2222
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8391.6666
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5841.6666
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 34597.406
$ws.Range("I134").Value = 3237.2666
$ws.Range("K134").Value = 9711.799800000001
$ws.Range("M134").Value = -7176.799800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4644.9473
$ws.Range("I99").Value = 3317.7778
$ws.Range("K99").Value = 3317.7778
$ws.Range("M99").Value = -1819.7778
$ws.Range("H126").Value = 4644.9473
$ws.Range("I126").Value = 3317.7778
$ws.Range("K126").Value = 9953.3334
$ws.Range("M126").Value = -7483.3334
$ws.Range("H129").Value = 74937.5
$ws.Range("J129").Value = 74937.5
$ws.Range("L129").Value = 74937.5
$ws.Range("N129").Value = -84937.5
$ws.Range("H132").Value = 1806.6
$ws.Range("I132").Value = 1652.5294
$ws.Range("K132").Value = 4957.5882
$ws.Range("M132").Value = -2427.5882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 258
$ws.Range("I2").Value = 75.166664
$ws.Range("K2").Value = 450.999984
$ws.Range("M2").Value = -337.999984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 54645
$ws.Range("J32").Value = 54645
$ws.Range("L32").Value = 54645
$ws.Range("N32").Value = -55237
$ws.Range("H45").Value = 79957
$ws.Range("J45").Value = 79957
$ws.Range("L45").Value = 79957
$ws.Range("N45").Value = -81075
$ws.Range("H70").Value = 8513.666999999999
$ws.Range("I70").Value = 6466.3335
$ws.Range("K70").Value = 6466.3335
$ws.Range("M70").Value = -6196.3335
$ws.Range("H73").Value = 8513.666999999999
$ws.Range("I73").Value = 6466.3335
$ws.Range("K73").Value = 6466.3335
$ws.Range("M73").Value = -5530.3335
$ws.Range("H97").Value = 900.2273
$ws.Range("I97").Value = 928.1177
$ws.Range("K97").Value = 928.1177
$ws.Range("M97").Value = -432.1177
$ws.Range("H126").Value = 3352.6155
$ws.Range("I126").Value = 3092.6667
$ws.Range("K126").Value = 9278.000100000001
$ws.Range("M126").Value = -6808.000100000001
$ws.Range("H132").Value = 55940.523
$ws.Range("I132").Value = 8208.388999999999
$ws.Range("K132").Value = 24625.167
$ws.Range("M132").Value = -22095.167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5765.3
$ws.Range("I40").Value = 3540.6
$ws.Range("K40").Value = 3540.6
$ws.Range("M40").Value = -3404.6
$ws.Range("H46").Value = 2116.7
$ws.Range("I46").Value = 2257.4443
$ws.Range("K46").Value = 2257.4443
$ws.Range("M46").Value = -2069.4443
$ws.Range("H55").Value = 3604.2
$ws.Range("J55").Value = 5660.3335
$ws.Range("L55").Value = 5660.3335
$ws.Range("N55").Value = -6006.3335
$ws.Range("H61").Value = 4286.0454
$ws.Range("I61").Value = 4276.353
$ws.Range("J61").Value = 4319
$ws.Range("K61").Value = 4276.353
$ws.Range("L61").Value = 4319
$ws.Range("M61").Value = -4074.353
$ws.Range("N61").Value = -4723
$ws.Range("H93").Value = 111114510
$ws.Range("I93").Value = 200002260
$ws.Range("J93").Value = 4839.75
$ws.Range("K93").Value = 200002260
$ws.Range("L93").Value = 4839.75
$ws.Range("M93").Value = -200001012
$ws.Range("N93").Value = -7335.75
$ws.Range("H113").Value = 4286.0454
$ws.Range("I113").Value = 4276.353
$ws.Range("J113").Value = 4319
$ws.Range("K113").Value = 4276.353
$ws.Range("L113").Value = 4319
$ws.Range("M113").Value = -2106.353
$ws.Range("N113").Value = -8659
$ws.Range("H122").Value = 4325.647
$ws.Range("I122").Value = 3634.4
$ws.Range("K122").Value = 10903.2
$ws.Range("M122").Value = -8453.200000000001
$ws.Range("H132").Value = 7381.5
$ws.Range("I132").Value = 6962
$ws.Range("J132").Value = 11996
$ws.Range("K132").Value = 20886
$ws.Range("L132").Value = 35988
$ws.Range("M132").Value = -18356
$ws.Range("N132").Value = -41048
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null
$ws.Range("H135").Value = 73939.2
$ws.Range("J135").Value = 73939.2
$ws.Range("L135").Value = 73939.2
$ws.Range("N135").Value = -84079.2
$ws.Range("H136").Value = 838233.2
$ws.Range("I136").Value = 838233.2
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2514699.6
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2512149.6
$ws.Range("N136").Value = $null
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2792.2173
$ws.Range("I81").Value = 1374.6842
$ws.Range("K81").Value = 2749.3684
$ws.Range("M81").Value = -1688.3684
$ws.Range("H84").Value = 2792.2173
$ws.Range("I84").Value = 1374.6842
$ws.Range("K84").Value = 13746.842
$ws.Range("M84").Value = -8442.841999999999
$ws.Range("H126").Value = 1348.5
$ws.Range("I126").Value = 1134.6666
$ws.Range("K126").Value = 3403.9998
$ws.Range("M126").Value = -933.9998000000001
$ws.Range("H132").Value = 45378.207
$ws.Range("I132").Value = 2819.842
$ws.Range("K132").Value = 8459.526
$ws.Range("M132").Value = -5929.526
$ws.Range("H135").Value = 1000000000
$ws.Range("J135").Value = 1000000000
$ws.Range("L135").Value = 1000000000
$ws.Range("N135").Value = -1000010140
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null
